# Update cryptocurrency price/volume data (cryptos.xlsx) per latest scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.265.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.88%  "
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("E9").Value = "  +0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.844.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.79%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.629.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.281.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.11%  "
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("E21").Value = "  +1.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.90"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0519"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.49%  "
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.179.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.97%  "
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.802"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.791"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.756.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +13.09%  "
$ws.Range("E47").Value = "  +2.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.71"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.53%  "
$ws.Range("E51").Value = "  -0.36%  "
